$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9006006717681885
$ws.Range("B1").Value = 1.320859432220459
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.679266333580017
